$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.798.75"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "2.272.81"
$ws.Range("E3").Value = "  +4.68%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.93"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.81"
$ws.Range("E7").Value = "  +9.46%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.660"
$ws.Range("E9").Value = "  +17.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.69"
$ws.Range("E10").Value = "  +10.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.71"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0969"
$ws.Range("E12").Value = "  +4.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.42"
$ws.Range("E13").Value = "  +8.69%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "2.616.80"
$ws.Range("E15").Value = "  +5.00%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.886"
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.81"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").Value = "2.292.82"
$ws.Range("E18").Value = "  +5.73%  "
$ws.Range("D19").Value = "42.773.35"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("E20").Value = "  +7.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.12"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.59"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.08"
$ws.Range("E25").Value = "  +5.43%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.44"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.01"
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.49"
$ws.Range("E33").Value = "  +15.05%  "
$ws.Range("E34").Value = "  +5.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.63"
$ws.Range("E36").Value = "  +31.55%  "
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  +12.51%  "
$ws.Range("E39").Value = "  +5.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0315"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("E41").Value = "  +7.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.65"
$ws.Range("E42").Value = "  +14.78%  "
$ws.Range("E43").Value = "  +7.06%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.29"
$ws.Range("E44").Value = "  +10.07%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.206"
$ws.Range("E45").Value = "  +8.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.06"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.93"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("E51").Value = "  +4.81%  "
